$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the "Publikationsangaben" sheet (keeps all styles,
#    column widths, row heights, etc.) and place the copy right after
#    it, then rename it to "Erscheinungsverlauf".
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("Publikationsangaben")
$template.Copy($null, $template)
$new = $wb.Worksheets.Item("Publikationsangaben (2)")
$new.Name = "Erscheinungsverlauf"

# ------------------------------------------------------------------
# 2. Update the header/meta block on the new sheet.
# ------------------------------------------------------------------
$new.Range("B1").Value = "Erscheinungsverlauf"
$new.Range("B2").Value = "Anzeige des Erscheinungsverlaufs"
$new.Range("B3").Value = "Dates of Publication and/or Sequential Designation"
$new.Range("B4").Value = "Erscheinungsverlauf"
$new.Range("B5").Value = "Numbering"
$new.Range("B6").Value = "Ticket #131"

# Fix the "Ticket #xxx" hyperlink so it points at issue 131 instead of
# the copied issue 119 link.
$new.Hyperlinks.Delete()
$new.Hyperlinks.Add($new.Range("B6"), "http://redmine.thulb.uni-jena.de/issues/131", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Ticket #131")

# ------------------------------------------------------------------
# 3. Update the "Datenfelder" description row.
# ------------------------------------------------------------------
$new.Range("A9").Value = "262 `$a"
$new.Range("B9").Value = "Angaben zur Zählung von fortlaufenden Ressourcen (PICA 4025)"
$new.Range("C9").Value = ""

# ------------------------------------------------------------------
# 4. Update the sample/test row. Row 14 existed in the template
#    (second PPN example) but the new sheet only carries a single
#    example row, so row 14 is cleared entirely.
# ------------------------------------------------------------------
$new.Range("A13").Value = "262 `$a"
$new.Range("B13").Value = 233814418
$new.Range("C13").Value = "1.1996/97; 2.1999; 3.2001 - 10.2010[?]; auch mit durchgehender Nr.-Zählung"
$new.Range("A14:H14").Clear()

# Active cell / view state for the new sheet.
$new.Range("C14").Select()

# ------------------------------------------------------------------
# 5. Add the analogous new data row (15) to "Publikationsangaben"
#    itself, mirroring the style of the existing rows 13/14.
# ------------------------------------------------------------------
$pub = $wb.Worksheets.Item("Publikationsangaben")
$pub.Range("A14:F14").Copy()
$pub.Range("A15:F15").PasteSpecial(-4122)
$pub.Range("A15").Value = "264 `$a, `$b und `$c"
$pub.Range("B15").Value = 537824324
$pub.Range("C15").Value = "Langewiesen : Linus Wittich, 1994-`nLangewiesen : Inform-Verl, 1994-"
$pub.Range("B6").Select()
$pub.Range("B16").Select()

# ------------------------------------------------------------------
# 6. Misc. selection / scroll-position bookkeeping on the other
#    sheets, matching what a human editing session left behind.
# ------------------------------------------------------------------
$beschreibung = $wb.Worksheets.Item("Beschreibung")
$beschreibung.Range("A4").Select()
$beschreibung.Application.ActiveWindow.ScrollRow = 4
$beschreibung.Range("A5").Select()

$basis = $wb.Worksheets.Item("Basisklassifikation")
$basis.Range("A4").Select()
$basis.Application.ActiveWindow.ScrollRow = 4
$basis.Range("B2").Select()

$zitate = $wb.Worksheets.Item("Bibliographische Zitate")
$zitate.Range("A4").Select()
$zitate.Application.ActiveWindow.ScrollRow = 4
$zitate.Range("B4").Select()

$sprach = $wb.Worksheets.Item("Sprachangaben")
$sprach.Range("C13").Select()

# ------------------------------------------------------------------
# 7. Finally, make the new sheet the active / selected tab, which is
#    what the recorded workbook view (activeTab) reflects.
# ------------------------------------------------------------------
$new.Activate()
$new.Range("C14").Select()
